$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The tail of the date/quantity matrix in column A:B was truncated one row
# too early and the last row's date was wrong. Fix row 210 and restore the
# two rows that should follow it.

# Column A holds dates written as plain digit strings (e.g. "20210105"),
# not real Excel dates, so force text entry to avoid Excel re-interpreting
# them as numbers.
$datesRange = $ws.Range("A210:A212")
$datesRange.NumberFormat = "@"

$ws.Cells.Item(210, 1).Value = "20210127"
$ws.Cells.Item(210, 2).Value = 9260.0

$ws.Cells.Item(211, 1).Value = "20210218"
$ws.Cells.Item(211, 2).Value = 10990.0

$ws.Cells.Item(212, 1).Value = "20210325"
$ws.Cells.Item(212, 2).Value = 8013.0

# Restore General format now that the text values are committed, so we
# don't leave a stray "Text" number format applied to these cells.
$datesRange.NumberFormat = "General"
